$wb = $excel.ActiveWorkbook

# --- StatOutput sheet: add header row + data row (counts query result) ---
$statOutput = $wb.Worksheets.Item("StatOutput")

$statOutput.Range("A1").Value = "number_of_files"
$statOutput.Range("B1").Value = "number_of_sample"
$statOutput.Range("C1").Value = "number_of_cases"
$statOutput.Range("D1").Value = "number_of_study"

$statOutput.Range("A2").Value = "'331"
$statOutput.Range("B2").Value = "'136"
$statOutput.Range("C2").Value = "'144"
$statOutput.Range("D2").Value = "'2"

# --- StatOutput_Message sheet: append a second run's log block (rows 11-21) ---
$cypherQuery = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_designation IN ['COTC007B','NCATS-COP01']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

$statMessage = $wb.Worksheets.Item("StatOutput_Message")

$statMessage.Range("A11").Value = "Neo4j_URL:"
$statMessage.Range("A12").Value = "bolt://ncias-q2251-c.nci.nih.gov:7687"
$statMessage.Range("A13").Value = "User_name:"
$statMessage.Range("A14").Value = "neo4j"
$statMessage.Range("A15").Value = "PWD:"
$statMessage.Range("A16").Value = "icdcDBneo4j0"
$statMessage.Range("A17").Value = "Cypher:"
$statMessage.Range("A18").Value = $cypherQuery
$statMessage.Range("A19").Value = "Output:"
$statMessage.Range("A20").Value = "C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC03_Canine_Filter_Study-ALL_Neo4jData.xlsx"
$statMessage.Range("A21").Value = "C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC03_Canine_Filter_Study-ALL_Neo4jData.xlsx"
